$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 4 (old '3월 나들이 도시락 한 상' row). Rows 5-14 shift up to become rows 4-13.
$ws.Rows(4).Delete()

# Step 2: the old '제주항공...' and '남도장터...' rows (originally rows 12 and 13) are now
# at row 11. Deleting row 11 twice removes both of them, which shifts the old row 14
# ('그랜드 조선 부산 X SSG.COM X 또떠남') content up into row 11.
$ws.Rows(11).Delete()
$ws.Rows(11).Delete()

# Step 3: row 10 currently holds the old '아이오페 ...' event (originally row 11).
# Replace it entirely with the new '크리니크' event data.
$r10A = @'
크리니크 3/30(목) 8PM
'@
$ws.Cells.Item(10, 1).Value = $r10A

$r10B = @'
치크팝 단독특가 & 마스카라/수분크림 대용량 1+1 기획
'@
$ws.Cells.Item(10, 2).Value = $r10B

$r10C = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000005157&domainSiteNo=6005
'@
$ws.Cells.Item(10, 3).Value = $r10C

$r10D = @'
크리니크 @SSG.LIVE 3/30(목) 20:00PM
'@
$ws.Cells.Item(10, 4).Value = $r10D

$r10G = @'
['이벤트/쿠폰 > 크리니크 @SSG.LIVE 3/30(목) 20:00PM', '스마일클럽', '                     SSG.LIVE 사은품 지급 및 이벤트 혜택 당첨 주의사항', ' - 사은품 지급 및 이벤트 혜택 제공', ' - 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기']
'@
$ws.Cells.Item(10, 7).Value = $r10G

# Dates are set via a text-forced NumberFormat so Excel keeps them as literal text
# (matching the original sheet, which stores dates as plain text strings) instead of
# auto-converting them into date serial numbers.
$r10E = @'
2023-03-22
'@
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = $r10E
$ws.Cells.Item(10, 5).ClearFormats()

$r10F = @'
2023-03-30
'@
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = $r10F
$ws.Cells.Item(10, 6).ClearFormats()

# Step 4: row 11 now holds the old '그랜드 조선 부산 X SSG.COM X 또떠남' event data (originally
# row 14) after the shifts above. Only the headline (A) and subhead (B) text change;
# the link/title/dates/description (C-G) already carried over correctly.
$r11A = @'
그랜드조선부산 X 또떠남
'@
$ws.Cells.Item(11, 1).Value = $r11A

$r11B = @'
해운대 오션뷰 객실 단독특가 + 청구혜택까지
'@
$ws.Cells.Item(11, 2).Value = $r11B

